# 776-RBI-EI-DB-DL-REC-NON-RNI-CTRFD-DL-MD-TR-1-Late Repayment - Loan product
# - fix the product name label (missing hyphen after "776") on both the
#   input and output sheets
# - lower-case the "currency" row label
# - drop the trailing space from the "US Dollar" currency value and make it
#   use the same (green) value-cell style as the other account-mapping rows
# - leave the input sheet ("ProductLoanInput") selected/active, with the
#   currency row (A6:B6) highlighted

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

# Fix the "currency" label casing
$ws1.Range("A6").Value = "currency"

# Fix the product name (add missing hyphen) - appears on both sheets
$ws1.Range("B1").Value = "776-RBI-EI-DB-DL-REC-NON-RNI-CTRFD-DL-MD-TR-1-Late Repayment"
$ws2.Range("B1").Value = "776-RBI-EI-DB-DL-REC-NON-RNI-CTRFD-DL-MD-TR-1-Late Repayment"

# Fix the currency value - remove the trailing space and match the style
# used by the other account/value cells (e.g. B32 "Cash")
$ws1.Range("B32").Copy()
$ws1.Range("B6").PasteSpecial(-4122)  # xlPasteFormats
$ws1.Range("B6").Value = "US Dollar"

# Restore the input sheet as the active tab/selection
$ws1.Activate()
$ws1.Range("A6:B6").Select()
